$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (21, 25, 26) ---
$ws.Range("E21").Value = 8

$ws.Range("C25").Value = "El_NumFA"
$ws.Range("D25").Value = "TP"
$ws.Range("E25").Value = 5

$ws.Range("D26").Value = "TP"
$ws.Range("E26").Value = 5

# --- Add new rows 27-44 ---
# String cells written in the precise order the data was originally entered,
# so the shared-string table ends up in the same order.
$ws.Range("A27").Value = "2025-12-08T08:00"
$ws.Range("A30").Value = "2025-12-09T13:30"
$ws.Range("A28").Value = "2025-12-08T13:30"
$ws.Range("A29").Value = "2025-12-09T08:30"
$ws.Range("A31").Value = "2025-12-10T08:00"
$ws.Range("A32").Value = "2025-12-10T13:30"
$ws.Range("A33").Value = "2025-12-11T08:00"
$ws.Range("A34").Value = "2025-12-12T08:00"
$ws.Range("A35").Value = "2025-12-12T13:30"
$ws.Range("A36").Value = "2025-12-15T08:00"
$ws.Range("A37").Value = "2025-12-15T13:30"
$ws.Range("A38").Value = "2025-12-16T08:30"
$ws.Range("A39").Value = "2025-12-16T13:30"
$ws.Range("A40").Value = "2025-12-17T08:00"
$ws.Range("A41").Value = "2025-12-17T13:30"
$ws.Range("A42").Value = "2025-12-18T08:00"
$ws.Range("A43").Value = "2025-12-19T08:00"
$ws.Range("A44").Value = "2025-12-19T13:30"
$ws.Range("B27").Value = "2025-12-08T12:15"
$ws.Range("B28").Value = "2025-12-08T17:45"
$ws.Range("B29").Value = "2025-12-09T12:15"
$ws.Range("B30").Value = "2025-12-09T17:45"
$ws.Range("B31").Value = "2025-12-10T12:15"
$ws.Range("B32").Value = "2025-12-10T17:45"
$ws.Range("B33").Value = "2025-12-11T12:15"
$ws.Range("B34").Value = "2025-12-12T12:15"
$ws.Range("B35").Value = "2025-12-12T17:45"
$ws.Range("B36").Value = "2025-12-15T12:15"
$ws.Range("B37").Value = "2025-12-15T17:45"
$ws.Range("B38").Value = "2025-12-16T12:15"
$ws.Range("B39").Value = "2025-12-16T17:45"
$ws.Range("B40").Value = "2025-12-17T12:15"
$ws.Range("B41").Value = "2025-12-17T17:45"
$ws.Range("B42").Value = "2025-12-18T12:15"
$ws.Range("B43").Value = "2025-12-19T12:15"
$ws.Range("B44").Value = "2025-12-19T17:45"
$ws.Range("C27").Value = "Leçon com. AJ"
$ws.Range("C28").Value = "Ondel TR"
$ws.Range("C31").Value = "Exam EN Telecom"
$ws.Range("C32").Value = "Sys ML SM"
$ws.Range("C33").Value = "Exam EN Fond"
$ws.Range("C36").Value = "Auto JPO"
$ws.Range("C37").Value = "Leçon FA"
$ws.Range("C41").Value = "En Telecom JSM"
$ws.Range("D41").Value = "Visite"
$ws.Range("C42").Value = "TS/Ondel Exam CD"
$ws.Range("C43").Value = "Python JO"
$ws.Range("C44").Value = "Exam OS PV"

# Remaining string cells that reuse already-known shared strings
$ws.Range("D27").Value = "leçon"
$ws.Range("D28").Value = "TP"
$ws.Range("C29").Value = "GL MS/TR"
$ws.Range("D29").Value = "CM"
$ws.Range("C30").Value = "ElectroPhy SC"
$ws.Range("D30").Value = "CM"
$ws.Range("D31").Value = "Examen"
$ws.Range("D32").Value = "TP"
$ws.Range("D33").Value = "Examen"
$ws.Range("C34").Value = "El_NumFA"
$ws.Range("D34").Value = "TP"
$ws.Range("C35").Value = "OS PV"
$ws.Range("D35").Value = "TP"
$ws.Range("D36").Value = "TP"
$ws.Range("D37").Value = "leçon"
$ws.Range("C38").Value = "GL MS/TR"
$ws.Range("D38").Value = "CM"
$ws.Range("C39").Value = "ElectroPhy SC"
$ws.Range("D39").Value = "CM"
$ws.Range("C40").Value = "El_NumFA"
$ws.Range("D40").Value = "TP"
$ws.Range("D42").Value = "Examen"
$ws.Range("D43").Value = "TP"
$ws.Range("D44").Value = "Examen"

# Numeric cells (classroom_id, user_id, promo_id)
$ws.Range("E27").Value = 4
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("E30").Value = 8
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = 1
$ws.Range("E32").Value = 7
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 1
$ws.Range("E33").Value = 10
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 1
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 1
$ws.Range("E35").Value = 7
$ws.Range("F35").Value = 1
$ws.Range("G35").Value = 1
$ws.Range("E36").Value = 11
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 1
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 1
$ws.Range("E38").Value = 7
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 1
$ws.Range("E39").Value = 8
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 1
$ws.Range("E40").Value = 5
$ws.Range("F40").Value = 1
$ws.Range("G40").Value = 1
$ws.Range("E41").Value = 12
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 1
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 1
$ws.Range("E43").Value = 7
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 1
$ws.Range("E44").Value = 7
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 1

# --- View state: zoom + selection, matching the saved workbook/sheet view ---
$excel.ActiveWindow.Zoom = 86
$ws.Range("C44").Select()
